$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename existing sheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Entries Base").Name = "Entry Base"
$wb.Worksheets.Item("Tags").Name = "Tag Defs"

# ---------------------------------------------------------------------------
# 2. Update data on "Defs" sheet (rows 2 & 3)
# ---------------------------------------------------------------------------
$wsDefs = $wb.Worksheets.Item("Defs")

$wsDefs.Range("A2").Value = "lgr0q1t4-3lil"
$wsDefs.Range("B2").Value = "2023-04-21T15:42:45.553"
$wsDefs.Range("C2").Value = "lgr0q1te"
$wsDefs.Range("D2").Value = "'FALSE"
$wsDefs.Range("E2").Value = "05a8"
$wsDefs.Range("F2").Value = "afree"
$wsDefs.Range("G2").Value = "3️⃣"
$wsDefs.Range("H2").Value = "Edited with description!"
$wsDefs.Range("I2").Value = "SECOND"

$wsDefs.Range("A3").Value = "lgr0q1te-5odh"
$wsDefs.Range("B3").Value = "2023-04-21T15:42:45.554"
$wsDefs.Range("C3").Value = "lgr0q1te"
$wsDefs.Range("D3").Value = "'FALSE"
$wsDefs.Range("E3").Value = "7gor"
$wsDefs.Range("F3").Value = "Five"
$wsDefs.Range("G3").Value = "5️⃣"
$wsDefs.Range("H3").Value = "not in first file, added to second"
$wsDefs.Range("I3").Value = "SECOND"

# ---------------------------------------------------------------------------
# 3. Add new row on "Point Defs" sheet
# ---------------------------------------------------------------------------
$wsPointDefs = $wb.Worksheets.Item("Point Defs")

$wsPointDefs.Range("A2").Value = "lgr0q1te-9rqg"
$wsPointDefs.Range("B2").Value = "2023-04-21T15:42:45.555"
$wsPointDefs.Range("C2").Value = "lgr0q1tf"
$wsPointDefs.Range("D2").Value = $false
$wsPointDefs.Range("E2").Value = "e0bq"
$wsPointDefs.Range("F2").Value = "0pc6"
$wsPointDefs.Range("G2").Value = "updated label"
$wsPointDefs.Range("H2").Value = "☝️"
$wsPointDefs.Range("I2").Value = "Set a description"
$wsPointDefs.Range("J2").Value = "BOOL"
$wsPointDefs.Range("K2").Value = "COUNT"
$wsPointDefs.Range("L2").Value = "TEXT"

# ---------------------------------------------------------------------------
# 4. Extend "Entry Base" header with new columns (_eid, _period before _note)
# ---------------------------------------------------------------------------
$wsEntryBase = $wb.Worksheets.Item("Entry Base")

$wsEntryBase.Range("F1").Value = "_eid"
$wsEntryBase.Range("G1").Value = "_period"
$wsEntryBase.Range("H1").Value = "_note"

# ---------------------------------------------------------------------------
# 5. Extend "Entry Points" header with new column (_eid before _val)
# ---------------------------------------------------------------------------
$wsEntryPoints = $wb.Worksheets.Item("Entry Points")

$wsEntryPoints.Range("G1").Value = "_eid"
$wsEntryPoints.Range("H1").Value = "_val"

# ---------------------------------------------------------------------------
# 6. Extend "Tag Defs" header with new columns (_emoji, _desc)
# ---------------------------------------------------------------------------
$wsTagDefs = $wb.Worksheets.Item("Tag Defs")

$wsTagDefs.Range("G1").Value = "_emoji"
$wsTagDefs.Range("H1").Value = "_desc"

# ---------------------------------------------------------------------------
# 7. Add new "Tags" sheet at the end
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNewTags = $wb.Worksheets.Add($null, $lastSheet)
$wsNewTags.Name = "Tags"

$wsNewTags.Range("A1").Value = "_uid"
$wsNewTags.Range("B1").Value = "_created"
$wsNewTags.Range("C1").Value = "_updated"
$wsNewTags.Range("D1").Value = "_deleted"
$wsNewTags.Range("E1").Value = "_did"
$wsNewTags.Range("F1").Value = "_pid"
$wsNewTags.Range("G1").Value = "tid"
